# Adjust move animation and balance, adjust blockspark spawn position
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# "Dashing Left Punch" raw stats live on row 35 (mirrored/read by row 11 via
# shared formulas). Bump the impact frame (B) and length (E); onhit/onblock
# (F/G) are formula-driven off these and recalc automatically.
$ws.Range("B35").Value = 13
$ws.Range("E35").Value = 30

# "Backdashing Left Kick" raw stats live on row 44 (mirrored by row 20).
# Adjust onhit (C) which drives the blockspark-related onblock/length calcs.
$ws.Range("C44").Value = 22

# Move the active selection to C45 to match where editing finished.
$ws.Range("C45").Select()
